# rebuilt tags model (v2.1.0)
$wb = $excel.ActiveWorkbook

$wsPackages   = $wb.Worksheets.Item("packages")
$wsAttributes = $wb.Worksheets.Item("attributes")
$wsTags       = $wb.Worksheets.Item("tags")

function Remove-HyperlinkAt($ws, $targetAddr) {
    # Deleting from a live COM collection while enumerating it is unsafe, so
    # repeatedly rescan-and-delete-one until no more matches are found.
    $found = $true
    while ($found) {
        $found = $false
        foreach ($h in @($ws.Hyperlinks)) {
            if ($h.Range.Address() -eq $targetAddr) {
                $h.Delete()
                $found = $true
                break
            }
        }
    }
}

# 1. Bump the package description / version string
$wsPackages.Range("C2").Value = "Mapping tables for processing raw data into unified model terminology (v2.1.0, 2022-06-29)"

# 2. attributes sheet: the "tags" column (G) used to hyperlink straight to the
#    ontology term IRI. The tags model was rebuilt so these references are no
#    longer rendered as hyperlinks - drop the link + hyperlink styling but
#    keep the text value as-is.
$attributeTagCells = '$G$2', '$G$3', '$G$4'
foreach ($addr in $attributeTagCells) {
    Remove-HyperlinkAt $wsAttributes $addr
}
$wsAttributes.Range("G2:G4").Style = "Normal"

# 3. tags sheet: identifiers switch from the full IRI (e.g.
#    "http://purl.obolibrary.org/obo/NCIT_C25415") / colon-form code
#    ("NCIT:C25415") to the short underscore-form code ("NCIT_C25415").
#    The identifier column (A) no longer hyperlinks out (objectIRI column C
#    keeps that job), and the label column (B) mirrors the new identifier.
$tagRows = @(
    @{ Row = 2; Code = "NCIT_C25415" },
    @{ Row = 3; Code = "NCIT_C25516" },
    @{ Row = 4; Code = "NCIT_C65107" }
)

foreach ($entry in $tagRows) {
    $r = $entry.Row
    $code = $entry.Code

    $idAddr = '$A$' + $r
    Remove-HyperlinkAt $wsTags $idAddr

    $idCell = $wsTags.Range("A$r")
    $idCell.Value = $code
    $idCell.Style = "Normal"

    $wsTags.Range("B$r").Value = $code
}
